# chore: update Sheets via scheduled runner
#
# Refreshes the market-board derived columns (H:N = currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) for a handful
# of leve rows across the per-job sheets, as produced by the scheduled
# price-refresh runner.

$wb = $excel.ActiveWorkbook

# --- ALC ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H111").Value = 1361.2941
$ws.Range("I111").Value = 1165.3334
$ws.Range("J111").Value = 1581.75
$ws.Range("K111").Value = 3496.0002
$ws.Range("L111").Value = 4745.25
$ws.Range("M111").Value = -429.0001999999999
$ws.Range("N111").Value = -10879.25

$ws.Range("H113").Value = 3082.2
$ws.Range("I113").Value = 2185.7144
$ws.Range("J113").Value = 5174
$ws.Range("K113").Value = 2185.7144
$ws.Range("L113").Value = 5174
$ws.Range("M113").Value = 1068.2856
$ws.Range("N113").Value = -11682

$ws.Range("H132").Value = 2912972.2
$ws.Range("I132").Value = 3478405.8
$ws.Range("J132").Value = 5028.5713
$ws.Range("K132").Value = 10435217.4
$ws.Range("L132").Value = 15085.7139
$ws.Range("M132").Value = -10432687.4
$ws.Range("N132").Value = -20145.7139

$ws.Range("H137").Value = 2079.3225
$ws.Range("I137").Value = 2138.9592
$ws.Range("J137").Value = 1854.5385
$ws.Range("K137").Value = 6416.8776
$ws.Range("L137").Value = 5563.6155
$ws.Range("M137").Value = -3866.8776
$ws.Range("N137").Value = -10663.6155

$ws.Range("H138").Value = 2388.1636
$ws.Range("I138").Value = 1169.8928
$ws.Range("K138").Value = 3509.6784
$ws.Range("M138").Value = 1630.3216

# --- ARM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H45").Value = 2866.6667
$ws.Range("I45").Value = 2866.6667
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2866.6667
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2489.6667
$ws.Range("N45").ClearContents()

$ws.Range("H61").Value = 6704.6113
$ws.Range("I61").Value = 6980.647
$ws.Range("K61").Value = 6980.647
$ws.Range("M61").Value = -6768.647

$ws.Range("H110").Value = 896.7778
$ws.Range("I110").Value = 896.7778
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 896.7778
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1148.2222
$ws.Range("N110").ClearContents()

$ws.Range("H132").Value = 4513.768
$ws.Range("I132").Value = 3161.1316
$ws.Range("J132").Value = 7369.3335
$ws.Range("K132").Value = 9483.3948
$ws.Range("L132").Value = 22108.0005
$ws.Range("M132").Value = -6953.3948
$ws.Range("N132").Value = -27168.0005

$ws.Range("H136").Value = 6704.6113
$ws.Range("I136").Value = 6980.647
$ws.Range("K136").Value = 20941.941
$ws.Range("M136").Value = -18391.941

# --- CRP -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 26840.25
$ws.Range("I31").Value = 51177.91
$ws.Range("J31").Value = 2502.5908
$ws.Range("K31").Value = 51177.91
$ws.Range("L31").Value = 2502.5908
$ws.Range("M31").Value = -50882.91
$ws.Range("N31").Value = -3092.5908

$ws.Range("H34").Value = 26840.25
$ws.Range("I34").Value = 51177.91
$ws.Range("J34").Value = 2502.5908
$ws.Range("K34").Value = 51177.91
$ws.Range("L34").Value = 2502.5908
$ws.Range("M34").Value = -50975.91
$ws.Range("N34").Value = -2906.5908

$ws.Range("H86").Value = 2445.9565
$ws.Range("I86").Value = 1937.75
$ws.Range("J86").Value = 3000.3635
$ws.Range("K86").Value = 1937.75
$ws.Range("L86").Value = 3000.3635
$ws.Range("M86").Value = -814.75
$ws.Range("N86").Value = -5246.363499999999

$ws.Range("H89").Value = 2445.9565
$ws.Range("I89").Value = 1937.75
$ws.Range("J89").Value = 3000.3635
$ws.Range("K89").Value = 9688.75
$ws.Range("L89").Value = 15001.8175
$ws.Range("M89").Value = -4072.75
$ws.Range("N89").Value = -26233.8175

$ws.Range("H111").Value = 40702
$ws.Range("J111").Value = 40702
$ws.Range("L111").Value = 40702
$ws.Range("N111").Value = -48882

$ws.Range("H134").Value = 7132.647
$ws.Range("I134").Value = 4958.815
$ws.Range("J134").Value = 15517.429
$ws.Range("K134").Value = 14876.445
$ws.Range("L134").Value = 46552.287
$ws.Range("M134").Value = -12341.445
$ws.Range("N134").Value = -51622.287

# --- CUL -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H5").Value = 854.3555
$ws.Range("I5").Value = 534.8570999999999
$ws.Range("J5").Value = 1133.9166
$ws.Range("K5").Value = 1604.5713
$ws.Range("L5").Value = 3401.7498
$ws.Range("M5").Value = -1492.5713
$ws.Range("N5").Value = -3625.7498

$ws.Range("H118").Value = 2388.65
$ws.Range("I118").Value = 724.5
$ws.Range("J118").Value = 2573.5557
$ws.Range("K118").Value = 2173.5
$ws.Range("L118").Value = 7720.6671
$ws.Range("M118").Value = -930.5
$ws.Range("N118").Value = -10206.6671

$ws.Range("H122").Value = 1845.6666
$ws.Range("I122").Value = 2125.6667
$ws.Range("J122").Value = 1565.6666
$ws.Range("K122").Value = 19131.0003
$ws.Range("L122").Value = 14090.9994
$ws.Range("M122").Value = -16681.0003
$ws.Range("N122").Value = -18990.9994

$ws.Range("H135").Value = 854.3555
$ws.Range("I135").Value = 534.8570999999999
$ws.Range("J135").Value = 1133.9166
$ws.Range("K135").Value = 4813.7139
$ws.Range("L135").Value = 10205.2494
$ws.Range("M135").Value = -2278.7139
$ws.Range("N135").Value = -15275.2494

# --- GSM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H2").Value = 63.42857
$ws.Range("I2").Value = 31.444445
$ws.Range("J2").Value = 121
$ws.Range("K2").Value = 31.444445
$ws.Range("L2").Value = 121
$ws.Range("M2").Value = 81.555555
$ws.Range("N2").Value = -347

$ws.Range("H102").Value = 1558.1364
$ws.Range("I102").Value = 1299.2858
$ws.Range("J102").Value = 2011.125
$ws.Range("K102").Value = 1299.2858
$ws.Range("L102").Value = 2011.125
$ws.Range("M102").Value = 322.7141999999999
$ws.Range("N102").Value = -5255.125

$ws.Range("H113").Value = 700
$ws.Range("I113").Value = 700
$ws.Range("K113").Value = 700
$ws.Range("M113").Value = 1470

$ws.Range("H133").Value = 63035
$ws.Range("J133").Value = 63035
$ws.Range("L133").Value = 63035
$ws.Range("N133").Value = -73155

# --- LTW -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H36").Value = 40532
$ws.Range("J36").Value = 40532
$ws.Range("L36").Value = 40532
$ws.Range("N36").Value = -41656

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46:N46").ClearContents()

# --- WVR -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H16").Value = 49999.75
$ws.Range("J16").Value = 49999.75
$ws.Range("L16").Value = 49999.75
$ws.Range("N16").Value = -50583.75

$ws.Range("H107").Value = 419.65
$ws.Range("I107").Value = 374.5
$ws.Range("J107").Value = 464.8
$ws.Range("K107").Value = 1123.5
$ws.Range("L107").Value = 1394.4
$ws.Range("M107").Value = 796.5
$ws.Range("N107").Value = -5234.4
